$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.424501
$ws.Range("H2").Value = 10.273503
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.408252
$ws.Range("N2").Value = 19.224756
$ws.Range("O2").Value = 0.8583439096634812
$ws.Range("P2").Value = 0.8583439096634812
$ws.Range("Q2").Value = 21.945065382252
$ws.Range("R2").Value = 197.505588440268
$ws.Range("S2").Value = 0.8583439096634812
$ws.Range("T2").Value = 0.8583439096634812

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.424501
$ws.Range("H3").Value = 10.273503
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6824433333333334
$ws.Range("N3").Value = 2.04733
$ws.Range("O3").Value = 0.09140887075868921
$ws.Range("P3").Value = 0.09140887075868921
$ws.Range("Q3").Value = 2.337027877443334
$ws.Range("R3").Value = 21.03325089699
$ws.Range("S3").Value = 0.09140887075868921
$ws.Range("T3").Value = 0.09140887075868921

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.424501
$ws.Range("H4").Value = 10.273503
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3751373333333333
$ws.Range("N4").Value = 1.125412
$ws.Range("O4").Value = 0.05024721957782962
$ws.Range("P4").Value = 0.05024721957782963
$ws.Range("Q4").Value = 1.284658173137333
$ws.Range("R4").Value = 11.561923558236
$ws.Range("S4").Value = 0.05024721957782962
$ws.Range("T4").Value = 0.05024721957782963
